# Generate Report for Handoff
#
# The localization status flips from "In Translation" to "Ready for
# handoff" and the "Latest Handoff Datetime" / "Latest HO Xliff Generate
# Date" timestamps advance ~1 minute, across the Overview roll-up sheet
# and the per-locale (zh-cn / de-de) detail sheets. Widening the status
# text also pushes Excel to re-autofit the narrower Status-ish columns.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn (E) / de-de (F) status columns, and the
# "Latest HO Xliff Generate Date" column (G).
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-13 10:48:04"

# zh-cn sheet: Status (C) and Latest Handoff Datetime (H)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-13 10:47:55"

# de-de sheet: Status (C) and Latest Handoff Datetime (H)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-13 10:48:04"

# The new, longer status text ("Ready for handoff") no longer fits the
# old column width, so re-fit those columns.
$overview.Columns("E:F").ColumnWidth = 16.33
$zhcn.Columns("C:C").ColumnWidth = 16.33
$dede.Columns("C:C").ColumnWidth = 16.33
